$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2:K2").ClearContents()
$ws.Range("I3:K3").ClearContents()
$ws.Range("K5").ClearContents()
$ws.Range("J6:K6").ClearContents()
$ws.Range("I7:K7").ClearContents()
